$d = $word.ActiveDocument

# Helper: force a run boundary at a given character position without
# changing any text, by dropping a temporary bookmark there and removing
# it again. (Adding the bookmark splits whatever run currently spans that
# position into two; deleting the bookmark does not re-merge them.)
function Split-At {
    param($pos)
    $bmName = "zzz_split_" + $pos
    $d.Bookmarks.Add($bmName, $d.Range($pos, $pos)) | Out-Null
    $d.Bookmarks($bmName).Delete()
}

# ---------------------------------------------------------------------
# 1) Title: "Definition" -> "Description"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Definition", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Description", 2) | Out-Null

# The replace above can coalesce "Data Category Module " and "Description"
# into one run; re-split them so the run layout matches the original
# (two runs, as before the edit).
$p1 = $d.Paragraphs.Item(1)
$t1 = $p1.Range.Text
$descIdx = $t1.IndexOf("Description")
if ($descIdx -gt 0) {
    Split-At ($p1.Range.Start + $descIdx)
}

# ---------------------------------------------------------------------
# 2) Namespace line: split "urn:iso:std:iso:30042:ed:3.0" into
#    "urn:" | "iso:std" | ":iso:30042:ed:3.0"  (text itself is unchanged)
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$t4 = $p4.Range.Text
$urnIdx = $t4.IndexOf("urn:iso:std:iso:30042:ed:3.0")
if ($urnIdx -ge 0) {
    $base4 = $p4.Range.Start + $urnIdx
    # split right-to-left so each newly-created run keeps a clean xml:space
    Split-At ($base4 + 11)
    Split-At ($base4 + 4)
}

# ---------------------------------------------------------------------
# 3) Description paragraph: re-home the "_GoBack" bookmark into the middle
#    of the sentence (right after "...valid TBX") and split the sentence
#    into its logical runs (text itself is unchanged).
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(6)
$t6 = $p6.Range.Text
$base6 = $p6.Range.Start

$quoteChar = [char]0x201D
$needle = "=" + $quoteChar + "term"
$needleIdx = $t6.IndexOf($needle)
$typeIdx = $t6.IndexOf("type=")
$dialectsIdx = $t6.IndexOf(" dialects")

if ($needleIdx -ge 0) {
    Split-At ($base6 + $needleIdx + $needle.Length)
}
if ($typeIdx -ge 0) {
    Split-At ($base6 + $typeIdx + 4)
}
if ($dialectsIdx -ge 0) {
    # Re-adding a bookmark named "_GoBack" moves it here (bookmark names
    # are unique), which also removes it from its old location after
    # "...or ISO 30042" and splits the run at this position.
    $goBackPos = $base6 + $dialectsIdx
    $d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos)) | Out-Null
}
